# Update cell values in Sheet1 to reflect the new TPM-derived calculations
# for the Ccl25-Ccr10 ligand-receptor pair (rows 2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (ECs -> FAPs)
$ws.Range("G2").Value = 2.738607
$ws.Range("H2").Value = 8.215821
$ws.Range("I2").Value = 0.2235648590725649
$ws.Range("J2").Value = 0.223564859072565
$ws.Range("P2").Value = 0.3930722505137151
$ws.Range("Q2").Value = 2.838827236034
$ws.Range("R2").Value = 25.549445124306
$ws.Range("S2").Value = 0.08787714229143463
$ws.Range("T2").Value = 0.08787714229143467

# Row 3 (ECs -> MuSCs)
$ws.Range("G3").Value = 2.738607
$ws.Range("H3").Value = 8.215821
$ws.Range("I3").Value = 0.2235648590725649
$ws.Range("J3").Value = 0.223564859072565
$ws.Range("P3").Value = 0.6069277494862849
$ws.Range("Q3").Value = 4.383323990169
$ws.Range("R3").Value = 39.449915911521
$ws.Range("S3").Value = 0.1356877167811303
$ws.Range("T3").Value = 0.1356877167811303

# Row 4 (FAPs -> FAPs)
$ws.Range("I4").Value = 0.5119261090069511
$ws.Range("J4").Value = 0.5119261090069511
$ws.Range("P4").Value = 0.3930722505137151
$ws.Range("R4").Value = 58.50395309903
$ws.Range("S4").Value = 0.2012239477640917
$ws.Range("T4").Value = 0.2012239477640917

# Row 5 (FAPs -> MuSCs)
$ws.Range("I5").Value = 0.5119261090069511
$ws.Range("J5").Value = 0.5119261090069511
$ws.Range("P5").Value = 0.6069277494862849
$ws.Range("R5").Value = 90.33370466635499
$ws.Range("S5").Value = 0.3107021612428594
$ws.Range("T5").Value = 0.3107021612428594

# Row 6 (MuSCs -> FAPs)
$ws.Range("I6").Value = 0.2645090319204839
$ws.Range("J6").Value = 0.2645090319204839
$ws.Range("P6").Value = 0.3930722505137151
$ws.Range("S6").Value = 0.1039711604581887
$ws.Range("T6").Value = 0.1039711604581887

# Row 7 (MuSCs -> MuSCs)
$ws.Range("I7").Value = 0.2645090319204839
$ws.Range("J7").Value = 0.2645090319204839
$ws.Range("P7").Value = 0.6069277494862849
$ws.Range("Q7").Value = 5.186095838331666
$ws.Range("S7").Value = 0.1605378714622952
$ws.Range("T7").Value = 0.1605378714622952
